# This script reproduces the crypto-price/volume refresh captured by the commit
# "Updated cryptos list ... with GitHub Actions": it rewrites the Price (D) and
# Volume(1h) (E) columns for the existing coin rows, and fixes the NEARProtocol /
# EthereumClassic rows (31-32), which had their data swapped upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain number (e.g. "0.451").
# Force their number format to Text ("@") *before* assigning the value so Excel
# keeps storing them as text (matching the original inline-string/text cells)
# instead of silently converting them into numeric cells.
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).NumberFormat = "@"

# --- Apply the updated values row by row ---
# Row 2
$ws.Cells.Item(2, 4).Value = '67.403.40'
$ws.Cells.Item(2, 5).Value = '  -1.39%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.770.07'
$ws.Cells.Item(3, 5).Value = '  -0.13%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.10%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '593.55'
$ws.Cells.Item(5, 5).Value = '  -0.42%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '165.90'
$ws.Cells.Item(6, 5).Value = '  -1.48%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '3.767.85'
$ws.Cells.Item(7, 5).Value = '  -0.11%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.02%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -0.12%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -0.91%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -1.43%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '0.451'
$ws.Cells.Item(12, 5).Value = '  +0.55%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '0.0000257'
$ws.Cells.Item(13, 5).Value = '  -2.87%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '35.87'
$ws.Cells.Item(14, 5).Value = '  -1.57%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '4.402.30'
$ws.Cells.Item(15, 5).Value = '  -0.17%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '3.785.23'
$ws.Cells.Item(16, 5).Value = '  +0.16%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '67.488.29'
$ws.Cells.Item(17, 5).Value = '  -1.21%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '17.74'
$ws.Cells.Item(18, 5).Value = '  -2.60%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.30%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '6.94'
$ws.Cells.Item(20, 5).Value = '  -1.53%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '10.56'
$ws.Cells.Item(21, 5).Value = '  -3.93%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '458.08'
$ws.Cells.Item(22, 5).Value = '  -2.21%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '0.696'

# Row 24
$ws.Cells.Item(24, 5).Value = '  +6.53%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '83.33'
$ws.Cells.Item(25, 5).Value = '  -1.76%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -4.33%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '11.81'
$ws.Cells.Item(27, 5).Value = '  -3.09%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -1.89%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -0.02%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -1.13%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'EthereumClassic'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(31, 4).Value = '29.76'
$ws.Cells.Item(31, 5).Value = '  -1.04%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'NEARProtocol'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(32, 4).Value = '7.21'
$ws.Cells.Item(32, 5).Value = '  -2.64%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -2.54%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '9.17'
$ws.Cells.Item(34, 5).Value = '  -1.04%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '0.999'
$ws.Cells.Item(35, 5).Value = '  -0.03%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '3.721.57'
$ws.Cells.Item(36, 5).Value = '  -0.20%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '0.100'
$ws.Cells.Item(37, 5).Value = '  -1.53%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -1.01%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  -1.07%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '0.996'
$ws.Cells.Item(40, 5).Value = '  -0.38%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '5.74'
$ws.Cells.Item(41, 5).Value = '  -1.27%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  +0.06%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -0.01%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '44.14'
$ws.Cells.Item(44, 5).Value = '  +0.51%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -2.70%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '46.88'
$ws.Cells.Item(46, 5).Value = '  +2.89%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -3.06%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '8.36'
$ws.Cells.Item(48, 5).Value = '  -2.60%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '145.88'
$ws.Cells.Item(49, 5).Value = '  +0.07%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '391.84'
$ws.Cells.Item(50, 5).Value = '  -3.66%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '2.756.66'
$ws.Cells.Item(51, 5).Value = '  +2.85%  '
